$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (class "1")
$ws.Range("B2").Value = 0.8118518518518518
$ws.Range("C2").Value = 0.9489177489177489
$ws.Range("D2").Value = 0.8750499001996008
$ws.Range("E2").Value = 1155

# Row 3 (class "2")
$ws.Range("B3").Value = 0.9362363919129082
$ws.Range("C3").Value = 0.9304482225656878
$ws.Range("D3").Value = 0.9333333333333333

# Row 4 (class "3")
$ws.Range("B4").Value = 0.8391959798994975
$ws.Range("C4").Value = 0.8487928843710292
$ws.Range("D4").Value = 0.8439671509791535
$ws.Range("E4").Value = 787

# Row 5 (class "4")
$ws.Range("B5").Value = 0.9139072847682119
$ws.Range("C5").Value = 0.3931623931623932
$ws.Range("D5").Value = 0.5498007968127491
$ws.Range("E5").Value = 351

# Row 6 (accuracy)
$ws.Range("B6").Value = 0.8517006802721089
$ws.Range("C6").Value = 0.8517006802721089
$ws.Range("D6").Value = 0.8517006802721089
$ws.Range("E6").Value = 0.8517006802721089

# Row 7 (macro avg)
$ws.Range("B7").Value = 0.8752978771081172
$ws.Range("C7").Value = 0.7803303122542148
$ws.Range("D7").Value = 0.8005377953312092

# Row 8 (weighted avg)
$ws.Range("B8").Value = 0.8587287508813221
$ws.Range("C8").Value = 0.8517006802721089
$ws.Range("D8").Value = 0.8407250098296173
